$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen), P (Precio $/Kg).
# These reflect a reshuffle of the existing 20 data rows (rows 2-21).

$rowData = @{
  2  = @{ D = 44875; J = 90;  K = 7000; L = 7000; M = 7000; O = "Provincia de Quillota"; P = 438 }
  3  = @{ D = 44208; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  4  = @{ D = 44882; J = 70;  K = 7000; L = 7000; M = 7000; O = "Provincia de Quillota"; P = 438 }
  5  = @{ D = 44188; J = 210; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  6  = @{ D = 44215; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  7  = @{ D = 44187; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  8  = @{ D = 44873; J = 250; K = 8000; L = 8000; M = 8000; O = "Provincia de Quillota"; P = 500 }
  9  = @{ D = 44186; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  10 = @{ D = 44210; J = 340; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  11 = @{ D = 44292; J = 90;  K = 6000; L = 6000; M = 6000; O = "Región Metropolitana";  P = 375 }
  12 = @{ D = 44189; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  13 = @{ D = 44204; J = 430; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  14 = @{ D = 44232; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  15 = @{ D = 44846; J = 250; K = 5000; L = 5000; M = 5000; O = "Provincia de Quillota"; P = 312 }
  16 = @{ D = 44231; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
  17 = @{ D = 44236; J = 180; K = 4000; L = 4500; M = 4167; O = "Región Metropolitana";  P = 260 }
  18 = @{ D = 44883; J = 180; K = 7000; L = 8000; M = 7500; O = "Provincia de Quillota"; P = 469 }
  19 = @{ D = 44855; J = 70;  K = 6000; L = 7000; M = 6500; O = "Provincia de Quillota"; P = 406 }
  20 = @{ D = 44251; J = 120; K = 5000; L = 5000; M = 5000; O = "Región Metropolitana";  P = 312 }
  21 = @{ D = 44230; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
}

foreach ($r in $rowData.Keys) {
  $row = $rowData[$r]
  $ws.Cells.Item($r, 4).Value  = $row.D   # D: Fecha
  $ws.Cells.Item($r, 10).Value = $row.J   # J: Volumen
  $ws.Cells.Item($r, 11).Value = $row.K   # K: Precio minimo
  $ws.Cells.Item($r, 12).Value = $row.L   # L: Precio maximo
  $ws.Cells.Item($r, 13).Value = $row.M   # M: Precio promedio ponderado
  $ws.Cells.Item($r, 15).Value = $row.O   # O: Origen
  $ws.Cells.Item($r, 16).Value = $row.P   # P: Precio $/Kg
}
